# Week 15 simulations update
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 (J.Allen)
$rushing.Range("C2").Value = 25
$rushing.Range("D2").Value = 27
$rushing.Range("E2").Value = 29
$rushing.Range("F2").Value = 23

# Row 3 (D.Singletary)
$rushing.Range("C3").Value = 59
$rushing.Range("E3").Value = 1

# Row 5 (M.Breida)
$rushing.Range("C5").Value = 12
$rushing.Range("D5").Value = 10
$rushing.Range("E5").Value = 3

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 (D.Singletary)
$receiving.Range("C2").Value = 36
$receiving.Range("D2").Value = 30
$receiving.Range("G2").Value = 4
$receiving.Range("H2").Value = 4

# Row 4 (M.Breida)
$receiving.Range("C4").Value = 8

# Row 6 (S.Diggs)
$receiving.Range("C6").Value = 93
$receiving.Range("D6").Value = 67
$receiving.Range("E6").Value = 29
$receiving.Range("F6").Value = 11
$receiving.Range("G6").Value = 21

# Row 7 (I.McKenzie)
$receiving.Range("E7").Value = 27
$receiving.Range("F7").Value = 12

# Row 8 (C.Beasley)
$receiving.Range("C8").Value = 85
$receiving.Range("D8").Value = 67

# Row 9 (G.Davis)
$receiving.Range("C9").Value = 23
$receiving.Range("D9").Value = 13
$receiving.Range("E9").Value = 14
$receiving.Range("F9").Value = 10
$receiving.Range("G9").Value = 11
$receiving.Range("H9").Value = 5

# Row 10 (J.Kumerow)
$receiving.Range("C10").Value = 7
$receiving.Range("D10").Value = 6

# Row 12 (D.Knox)
$receiving.Range("C12").Value = 41
$receiving.Range("D12").Value = 33
$receiving.Range("E12").Value = 15
$receiving.Range("F12").Value = 11
$receiving.Range("G12").Value = 13
$receiving.Range("H12").Value = 8
